$wb = $excel.ActiveWorkbook

# Rename the single worksheet from "FUNKY" to "FUNKY_DENOM"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "FUNKY_DENOM"
